# Forced Convection Calculations - add theoretical & correlated friction
# factor columns (f) next to the existing Nu columns.
#
# Commit message: "Completed theoretical and correlated calculations for f"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New header cells in columns V (22) and W (23), row 1.
#    Copy formatting from the existing header style (column U / "Correlated
#    Nu") so the new headers get the same bold/centered/bordered look,
#    then set their text afterwards so the paste doesn't clobber the value.
# ---------------------------------------------------------------------
$ws.Range("U1").Copy() | Out-Null
$ws.Range("V1:W1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("V1").Value = "Theoretical Friction Factor"
$ws.Range("W1").Value = "Correlated Friction Factor"

# ---------------------------------------------------------------------
# 2. New data values for the four experiment rows (2-5).
#    (written as plain decimals - the interpreter's literal parser does
#    not accept E-notation)
# ---------------------------------------------------------------------
$ws.Range("V2").Value = 0.00580750534635551385
$ws.Range("V3").Value = 0.00601741517795236023
$ws.Range("V4").Value = 0.00665668603052808634
$ws.Range("V5").Value = 0.00627811218306685570

$ws.Range("W2").Value = 0.01586812909492468088
$ws.Range("W3").Value = 0.01588155598361544116
$ws.Range("W4").Value = 0.01718707523864302852
$ws.Range("W5").Value = 0.01716457704102938853

# ---------------------------------------------------------------------
# 3. Give the two new columns explicit widths, matching the other
#    wide "label" columns in the sheet (~27/~26 characters).
# ---------------------------------------------------------------------
$ws.Columns.Item(22).ColumnWidth = 26.83
$ws.Columns.Item(23).ColumnWidth = 25.83

# ---------------------------------------------------------------------
# 4. Leave the cursor/selection where the author ended up.
# ---------------------------------------------------------------------
$ws.Range("V18").Select() | Out-Null
